$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.348.72"
$ws.Range("E2").Value = "  +0.46%  "
$ws.Range("D3").Value = "3.507.78"
$ws.Range("E3").Value = "  +0.04%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.95"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.52"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.74%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.487"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.34%  "
$ws.Range("E9").Value = "  +0.79%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.24"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.23%  "
$ws.Range("E11").Value = "  +1.64%  "
$ws.Range("D12").Value = "4.097.81"
$ws.Range("E12").Value = "  +0.43%  "
$ws.Range("E13").Value = "  +1.13%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000182"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.19%  "
$ws.Range("D15").Value = "3.503.34"
$ws.Range("E15").Value = "  +0.12%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.78"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -6.72%  "
$ws.Range("D17").Value = "64.343.09"
$ws.Range("E17").Value = "  +0.79%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "9.93"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.54%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.78"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.39%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.82"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "387.70"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.571"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.40%  "
$ws.Range("D23").Value = "3.644.63"
$ws.Range("E23").Value = "  +0.29%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "74.31"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.78%  "
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("E26").Value = "  +2.09%  "
$ws.Range("E27").Value = "  -0.32%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.48"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.46%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.997"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.14%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.33"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.13%  "
$ws.Range("E31").Value = "  -1.00%  "
$ws.Range("E32").Value = "  -5.28%  "
$ws.Range("D33").Value = "3.525.10"
$ws.Range("E33").Value = "  +0.60%  "
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("E35").Value = "  +2.94%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "23.58"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.92%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.24"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.52%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.57"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.96%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.91"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.97%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "163.35"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.39%  "
$ws.Range("E41").Value = "  -2.76%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.806"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.15%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "25.76"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.42%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.84"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.85%  "
$ws.Range("E46").Value = "  +0.03%  "
$ws.Range("B47").Value = "ONDO"
$ws.Range("C47").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.18"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.68%  "
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.66"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.07%  "
$ws.Range("D49").Value = "2.475.23"
$ws.Range("E49").Value = "  +1.29%  "
$ws.Range("E50").Value = "  -1.63%  "
$ws.Range("E51").Value = "  +0.30%  "
